# Generate Report for Handback
#
# Row 7 in both the "zh-cn" and "de-de" worksheets corresponds to the
# handoff file "1a50287c-500b-4403-bef3-d1bee53271ba". A handback has come
# in for it, but it turned out to be based on a stale handoff version, so:
#   - "Latest Target File"   (I7) gets the source .md file (hyperlinked,
#                              like every other populated row in col I/A)
#   - "Latest Handback File" (J7) gets the generated .xlf handback name
#   - "Latest Handback DateTime" (K7) gets the timestamp the handback
#                              came in
#   - "Error Detail"         (P7) gets a message explaining that the
#                              handback was based on an old handoff commit

$wb = $excel.ActiveWorkbook

$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/f864a7040b483fff59c5af2257a78e33b40f3acf/e2e/1a50287c-500b-4403-bef3-d1bee53271ba.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/30cc18902f82e39c55a671ffa764a857c8057bc7/e2e/1a50287c-500b-4403-bef3-d1bee53271ba.md."

# zh-cn -> 2016-08-12 03:16:26, de-de -> 2016-08-12 03:16:33
$sheetInfo = @{
    "zh-cn" = @{
        HandbackFile = "1a50287c-500b-4403-bef3-d1bee53271ba.7e793a050a5b4660054b4c760319ad1361ee7708.zh-cn.xlf"
        HandbackDate = "2016-08-12 03:16:26"
        TargetUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/30cc18902f82e39c55a671ffa764a857c8057bc7/e2e/1a50287c-500b-4403-bef3-d1bee53271ba.md"
    }
    "de-de" = @{
        HandbackFile = "1a50287c-500b-4403-bef3-d1bee53271ba.7e793a050a5b4660054b4c760319ad1361ee7708.de-de.xlf"
        HandbackDate = "2016-08-12 03:16:33"
        TargetUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/30cc18902f82e39c55a671ffa764a857c8057bc7/e2e/1a50287c-500b-4403-bef3-d1bee53271ba.md"
    }
}

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $info = $sheetInfo[$sheetName]

    $sourceFileName = $ws.Range("A7").Value2

    # "Latest Target File" mirrors the source file name, hyperlinked just
    # like column A and like every other already-processed row's column I.
    $ws.Range("I7").Value = $sourceFileName
    $ws.Range("I7").Font.Underline = $true
    $ws.Range("I7").Font.Color = 15570276

    $ws.Hyperlinks.Add($ws.Range("I7"), $info.TargetUrl, "", "", $sourceFileName) | Out-Null

    # "Latest Handback File"
    $ws.Range("J7").Value = $info.HandbackFile

    # "Latest Handback DateTime"
    $ws.Range("K7").Value = $info.HandbackDate

    # "Error Detail"
    $ws.Range("P7").Value = $errorMessage
}
